$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.062.99"
$ws.Range("E2").Value = "  -0.61%  "

# Row 3
$ws.Range("D3").Value = "2.374.26"
$ws.Range("E3").Value = "  -1.16%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").Value = "548.37"
$ws.Range("E5").Value = "  -0.45%  "

# Row 6
$ws.Range("D6").Value = "137.77"
$ws.Range("E6").Value = "  -3.45%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").Value = "0.526"
$ws.Range("E8").Value = "  -2.77%  "

# Row 9
$ws.Range("D9").Value = "2.375.81"
$ws.Range("E9").Value = "  -0.92%  "

# Row 10
$ws.Range("E10").Value = "  +1.02%  "

# Row 11
$ws.Range("E11").Value = "  +1.20%  "

# Row 12
$ws.Range("D12").Value = "5.36"
$ws.Range("E12").Value = "  +0.78%  "

# Row 13
$ws.Range("E13").Value = "  -0.45%  "

# Row 14
$ws.Range("D14").Value = "25.02"
$ws.Range("E14").Value = "  -2.40%  "

# Row 15
$ws.Range("D15").Value = "2.786.59"
$ws.Range("E15").Value = "  -1.60%  "

# Row 16
$ws.Range("D16").Value = "0.0000165"
$ws.Range("E16").Value = "  -0.98%  "

# Row 17
$ws.Range("D17").Value = "61.000.17"
$ws.Range("E17").Value = "  +0.38%  "

# Row 18
$ws.Range("D18").Value = "2.381.32"
$ws.Range("E18").Value = "  -0.70%  "

# Row 19
$ws.Range("D19").Value = "10.79"
$ws.Range("E19").Value = "  -0.56%  "

# Row 20
$ws.Range("D20").Value = "4.14"
$ws.Range("E20").Value = "  -0.47%  "

# Row 21
$ws.Range("D21").Value = "320.31"
$ws.Range("E21").Value = "  +0.08%  "

# Row 22
$ws.Range("D22").Value = "6.69"
$ws.Range("E22").Value = "  -1.00%  "

# Row 24
$ws.Range("D24").Value = "64.18"
$ws.Range("E24").Value = "  +0.56%  "

# Row 25
$ws.Range("E25").Value = "  -13.02%  "

# Row 26
$ws.Range("D26").Value = "8.54"
$ws.Range("E26").Value = "  +3.26%  "

# Row 27
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.02%  "

# Row 28
$ws.Range("D28").Value = "2.487.34"
$ws.Range("E28").Value = "  -1.17%  "

# Row 29
$ws.Range("D29").Value = "8.12"
$ws.Range("E29").Value = "  -0.51%  "

# Row 30
$ws.Range("D30").Value = "506.88"
$ws.Range("E30").Value = "  -5.94%  "

# Row 31
$ws.Range("E31").Value = "  +2.43%  "

# Row 32
$ws.Range("D32").Value = "0.0₃0882"
$ws.Range("E32").Value = "  -7.05%  "

# Row 33
$ws.Range("D33").Value = "1.38"
$ws.Range("E33").Value = "  -5.01%  "

# Row 34
$ws.Range("D34").Value = "1.82"
$ws.Range("E34").Value = "  -1.88%  "

# Row 35
$ws.Range("D35").Value = "1.51"
$ws.Range("E35").Value = "  -5.31%  "

# Row 36
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.12%  "

# Row 37
$ws.Range("D37").Value = "4.68"
$ws.Range("E37").Value = "  -1.71%  "

# Row 38
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "0.378"
$ws.Range("E38").Value = "  -0.27%  "

# Row 39
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "1.87"
$ws.Range("E39").Value = "  +0.36%  "

# Row 40
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").Value = "18.57"
$ws.Range("E40").Value = "  +2.06%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").Value = "5.34"
$ws.Range("E41").Value = "  -4.51%  "

# Row 42
$ws.Range("D42").Value = "145.81"
$ws.Range("E42").Value = "  +5.00%  "

# Row 43
$ws.Range("E43").Value = "  -0.07%  "

# Row 44
$ws.Range("D44").Value = "41.65"
$ws.Range("E44").Value = "  +3.34%  "

# Row 45
$ws.Range("D45").Value = "147.54"
$ws.Range("E45").Value = "  +3.80%  "

# Row 46
$ws.Range("D46").Value = "3.59"
$ws.Range("E46").Value = "  -1.59%  "

# Row 47
$ws.Range("D47").Value = "2.06"
$ws.Range("E47").Value = "  -7.01%  "

# Row 48
$ws.Range("D48").Value = "0.0521"
$ws.Range("E48").Value = "  -0.44%  "

# Row 49
$ws.Range("D49").Value = "19.27"
$ws.Range("E49").Value = "  -5.78%  "

# Row 50
$ws.Range("D50").Value = "0.575"
$ws.Range("E50").Value = "  -0.90%  "

# Row 51
$ws.Range("D51").Value = "0.0910"
$ws.Range("E51").Value = "  -0.17%  "
